$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) D15: "Specialty Chemicals" -> "Specialty Bust Chemicals"
$ws.Range("D15").Value = "Specialty Bust Chemicals"

# 2) B52: "Automated Data Processing" -> "Automatic Data Processing"
$ws.Range("B52").Value = "Automatic Data Processing"

# 3) Move the Fleetcor row (currently row 203) up to row 132 (alphabetical
#    re-sort after the company rebranded from Fleetcor/FLT to Corpay/CPAY),
#    shifting rows 132-202 down to 133-203.
$ws.Rows.Item(132).Insert()
$ws.Rows.Item(204).Cut($ws.Rows.Item(132))
$ws.Rows.Item(204).Delete()

# Rename the moved row to reflect the Corpay rebrand.
$ws.Range("A132").Value = "CPAY"
$ws.Range("B132").Value = "Corpay"
